# raley_brooks.xlsx regen: "K" column (col G, formerly "Strike#") is
# recalculated from the regenerated std/mean based s_vals pipeline.
# This writes the freshly-computed K values (and, for row 50, the
# corrected IP/I0/IF values) onto Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (column G) values, keyed by worksheet row number.
$kUpdates = @{
    2 = 2; 3 = 1; 4 = 0; 5 = 2; 6 = 1; 7 = 3;
    8 = 0; 9 = 2; 10 = 2; 11 = 2; 12 = 0; 13 = 0;
    14 = 0; 15 = 0; 16 = 1; 17 = 0; 18 = 0; 19 = 0;
    20 = 2; 21 = 0; 22 = 3; 23 = 1; 24 = 1; 25 = 0;
    26 = 2; 27 = 1; 28 = 0; 29 = 1; 30 = 0; 31 = 2;
    32 = 2; 33 = 2; 34 = 0; 35 = 1; 36 = 1; 37 = 1;
    38 = 0; 39 = 0; 40 = 0; 41 = 3; 42 = 0; 43 = 0;
    44 = 0; 45 = 2; 46 = 3; 47 = 1; 48 = 0; 49 = 2;
    51 = 2; 52 = 1; 53 = 0; 55 = 3; 56 = 2; 57 = 0;
    58 = 2; 60 = 0; 61 = 0; 62 = 0; 63 = 0; 64 = 1;
    65 = 0; 66 = 2; 67 = 1; 68 = 0
}

foreach ($row in $kUpdates.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $kUpdates[$row]
}

# Row 50 needed its K/IP recalculated, and its I0/IF (cols I/J) values
# were swapped as part of the regen.
$ws.Cells.Item(50, 7).Value = 0   # G50 (K)
$ws.Cells.Item(50, 8).Value = 3   # H50 (IP)
$ws.Cells.Item(50, 9).Value = 7   # I50 (I0)
$ws.Cells.Item(50, 10).Value = 9  # J50 (IF)
